$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Move the "Create task" row (old row 14) up into row 11, which currently
#    holds a stray duplicate of row 3 ("Get oauth access token..."). Also trim
#    the trailing whitespace on the JSON payload template text.
# ---------------------------------------------------------------------------
$newA11 = $ws.Range("A14").Value()
$newB11 = $ws.Range("B14").Value()
$newD11 = $ws.Range("D14").Value()
$newE11 = $ws.Range("E14").Value().TrimEnd()

$ws.Range("A11").Value = $newA11
$ws.Range("B11").Value = $newB11
$ws.Range("D11").Value = $newD11
$ws.Range("E11").Value = $newE11

$ws.Rows.Item(11).RowHeight = 43.2

# ---------------------------------------------------------------------------
# 2. Remove the now-duplicated row 14 entirely; this also shifts the trailing
#    spacer row (old row 17) up to row 16.
# ---------------------------------------------------------------------------
$ws.Rows.Item(14).Delete()

# ---------------------------------------------------------------------------
# 3. Rebuild the hyperlinks collection in the desired final order so the
#    relationship ids line up the same way the target workbook does.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B2"), "https://5fcc640551f70e00161f24b2.mockapi.io/api/v1/token")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://accounts.google.com/o/oauth2/token")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://v2.convertapi.com/upload")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://gorest.co.in/public-api/users")
$ws.Hyperlinks.Add($ws.Range("B9"), "https://gorest.co.in/public-api/users")
$ws.Hyperlinks.Add($ws.Range("B6:B8"), "https://gorest.co.in/public-api/users", $null, $null, "https://gorest.co.in/public-api/users")
$ws.Hyperlinks.Add($ws.Range("B12"), "https://maxsoft-mock-server-demo.web.app/photos")
$ws.Hyperlinks.Add($ws.Range("B13"), "https://api.apis.guru/", "version/#jsonFile")
$ws.Hyperlinks.Add($ws.Range("B11"), "http://8e4jz.mocklab.io/tasks")

# ---------------------------------------------------------------------------
# 4. Update the sheet selection to match the new state of the workbook.
# ---------------------------------------------------------------------------
$null = $ws.Range("B18").Select()
